$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The list of "Periodo Mora" values (column E, rows 16-22) is reversed
# from ascending (2305..2311) to descending (2311..2305), and the
# "Valor Mora" (column F) value that was attached to period 2311
# (43307) now travels with it to row 16, while row 22 reverts to the
# common value (46400) shared by the other rows.
$ws.Range("E16").Value = "2311"
$ws.Range("E17").Value = "2310"
$ws.Range("E18").Value = "2309"
$ws.Range("E19").Value = "2308"
$ws.Range("E20").Value = "2307"
$ws.Range("E21").Value = "2306"
$ws.Range("E22").Value = "2305"

$ws.Range("F16").Value = 43307
$ws.Range("F22").Value = 46400
